# Auto-generated script applying scheduled market-data refresh values
# to the Leve profit-tracking workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each touched row we update currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N)
# to the latest pulled market values. Some rows gain or lose a profit cell (M/N) depending on
# whether the NQ or HQ price is lower, matching the source data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 567.2222
$ws.Range("I2").Value = 386.42856
$ws.Range("K2").Value = 386.42856
$ws.Range("M2").Value = -273.42856
$ws.Range("H8").Value = 63.625
$ws.Range("I8").Value = 63.625
$ws.Range("K8").Value = 190.875
$ws.Range("M8").Value = -51.875
$ws.Range("H80").Value = 3525.25
$ws.Range("J80").Value = 5067.1665
$ws.Range("L80").Value = 15201.4995
$ws.Range("N80").Value = -17197.4995
$ws.Range("H83").Value = 3525.25
$ws.Range("J83").Value = 5067.1665
$ws.Range("L83").Value = 45604.4985
$ws.Range("N83").Value = -55588.4985
$ws.Range("H86").Value = 4500.5
$ws.Range("I86").Value = 3499.75
$ws.Range("J86").Value = 6502
$ws.Range("K86").Value = 3499.75
$ws.Range("L86").Value = 6502
$ws.Range("M86").Value = -2376.75
$ws.Range("N86").Value = -8748
$ws.Range("H89").Value = 4500.5
$ws.Range("I89").Value = 3499.75
$ws.Range("J89").Value = 6502
$ws.Range("K89").Value = 17498.75
$ws.Range("L89").Value = 32510
$ws.Range("M89").Value = -11882.75
$ws.Range("N89").Value = -43742
$ws.Range("H111").Value = 3550
$ws.Range("I111").Value = 3650
$ws.Range("K111").Value = 10950
$ws.Range("M111").Value = -7883
$ws.Range("H123").Value = 85000
$ws.Range("J123").Value = 85000
$ws.Range("L123").Value = 85000
$ws.Range("N123").Value = -94800
$ws.Range("H137").Value = 2518.2727
$ws.Range("J137").Value = 3999.8333
$ws.Range("L137").Value = 11999.4999
$ws.Range("N137").Value = -17099.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1749.1428
$ws.Range("J45").Value = 2000
$ws.Range("L45").Value = 2000
$ws.Range("N45").Value = -2754
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3629.9412
$ws.Range("I132").Value = 2644.111
$ws.Range("J132").Value = 4739
$ws.Range("K132").Value = 7932.333
$ws.Range("L132").Value = 14217
$ws.Range("M132").Value = -5402.333
$ws.Range("N132").Value = -19277

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 3500
$ws.Range("J15").Value = 3500
$ws.Range("L15").Value = 3500
$ws.Range("N15").Value = -3954
$ws.Range("H87").Value = 47545
$ws.Range("I87").Value = 47545
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 47545
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -46297
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 47545
$ws.Range("I90").Value = 47545
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 142635
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -136395
$ws.Range("N90").ClearContents()
$ws.Range("H134").Value = 2792.1
$ws.Range("I134").Value = 2654
$ws.Range("J134").Value = 3114.3333
$ws.Range("K134").Value = 7962
$ws.Range("L134").Value = 9342.999899999999
$ws.Range("M134").Value = -5427
$ws.Range("N134").Value = -14412.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1272.4546
$ws.Range("I16").Value = 1256.5555
$ws.Range("J16").Value = 1344
$ws.Range("K16").Value = 1256.5555
$ws.Range("L16").Value = 1344
$ws.Range("M16").Value = -969.5554999999999
$ws.Range("N16").Value = -1918
$ws.Range("H62").Value = 2641.6667
$ws.Range("I62").Value = 2641.6667
$ws.Range("K62").Value = 2641.6667
$ws.Range("M62").Value = -2017.6667
$ws.Range("H65").Value = 2641.6667
$ws.Range("I65").Value = 2641.6667
$ws.Range("K65").Value = 13208.3335
$ws.Range("M65").Value = -10088.3335
$ws.Range("H99").Value = 5725.4
$ws.Range("I99").Value = 5322.3335
$ws.Range("K99").Value = 5322.3335
$ws.Range("M99").Value = -3824.3335
$ws.Range("H113").Value = 1272.4546
$ws.Range("I113").Value = 1256.5555
$ws.Range("J113").Value = 1344
$ws.Range("K113").Value = 1256.5555
$ws.Range("L113").Value = 1344
$ws.Range("M113").Value = 913.4445000000001
$ws.Range("N113").Value = -5684
$ws.Range("H125").Value = 43830.25
$ws.Range("J125").Value = 43830.25
$ws.Range("L125").Value = 43830.25
$ws.Range("N125").Value = -48750.25
$ws.Range("H126").Value = 5725.4
$ws.Range("I126").Value = 5322.3335
$ws.Range("K126").Value = 15967.0005
$ws.Range("M126").Value = -13497.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 714817.1
$ws.Range("I4").Value = 833453.3
$ws.Range("K4").Value = 2500359.9
$ws.Range("M4").Value = -2500247.9
$ws.Range("H12").Value = 167250.61
$ws.Range("I12").Value = 375033
$ws.Range("J12").Value = 1024.7
$ws.Range("K12").Value = 1125099
$ws.Range("L12").Value = 3074.1
$ws.Range("M12").Value = -1124926
$ws.Range("N12").Value = -3420.1
$ws.Range("H39").Value = 3125.75
$ws.Range("J39").Value = 3125.75
$ws.Range("L39").Value = 9377.25
$ws.Range("N39").Value = -9965.25
$ws.Range("H50").Value = 284.75
$ws.Range("I50").Value = 284.75
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 854.25
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -373.25
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value = 284.75
$ws.Range("I53").Value = 284.75
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 854.25
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -373.25
$ws.Range("N53").ClearContents()
$ws.Range("H104").Value = 5029
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 5029
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 15087
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -20329
$ws.Range("H121").Value = 744.25
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H131").Value = 1599.9656
$ws.Range("I131").Value = 1079.6
$ws.Range("J131").Value = 1708.375
$ws.Range("K131").Value = 3238.8
$ws.Range("L131").Value = 5125.125
$ws.Range("M131").Value = 1801.2
$ws.Range("N131").Value = -15205.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()
$ws.Range("H80").Value = 2590.923
$ws.Range("I80").Value = 2525.5
$ws.Range("J80").Value = 2695.6
$ws.Range("K80").Value = 2525.5
$ws.Range("L80").Value = 2695.6
$ws.Range("M80").Value = -1527.5
$ws.Range("N80").Value = -4691.6
$ws.Range("H83").Value = 2590.923
$ws.Range("I83").Value = 2525.5
$ws.Range("J83").Value = 2695.6
$ws.Range("K83").Value = 12627.5
$ws.Range("L83").Value = 13478
$ws.Range("M83").Value = -7635.5
$ws.Range("N83").Value = -23462
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 10000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -10676
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 10000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -12340
$ws.Range("H132").Value = 6457.2104
$ws.Range("I132").Value = 5999.077
$ws.Range("K132").Value = 17997.231
$ws.Range("M132").Value = -15467.231
$ws.Range("H136").Value = 4894.75
$ws.Range("I136").Value = 4751.143
$ws.Range("J136").Value = 5900
$ws.Range("K136").Value = 14253.429
$ws.Range("L136").Value = 17700
$ws.Range("M136").Value = -11703.429
$ws.Range("N136").Value = -22800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29995
$ws.Range("I64").Value = 29985
$ws.Range("K64").Value = 29985
$ws.Range("M64").Value = -29737
$ws.Range("H67").Value = 29995
$ws.Range("I67").Value = 29985
$ws.Range("K67").Value = 29985
$ws.Range("M67").Value = -29127
$ws.Range("H70").Value = 24995
$ws.Range("I70").Value = 24990
$ws.Range("K70").Value = 24990
$ws.Range("M70").Value = -24675
$ws.Range("H73").Value = 24995
$ws.Range("I73").Value = 24990
$ws.Range("K73").Value = 24990
$ws.Range("M73").Value = -23898
$ws.Range("H75").Value = 24997.25
$ws.Range("I75").Value = 24994.5
$ws.Range("K75").Value = 24994.5
$ws.Range("M75").Value = -24058.5
$ws.Range("H78").Value = 24997.25
$ws.Range("I78").Value = 24994.5
$ws.Range("K78").Value = 74983.5
$ws.Range("M78").Value = -70303.5
$ws.Range("H86").Value = 49215
$ws.Range("J86").Value = 49215
$ws.Range("L86").Value = 49215
$ws.Range("N86").Value = -51461
$ws.Range("H89").Value = 49215
$ws.Range("J89").Value = 49215
$ws.Range("L89").Value = 246075
$ws.Range("N89").Value = -257307
$ws.Range("H132").Value = 3530.2222
$ws.Range("I132").Value = 3530.2222
$ws.Range("K132").Value = 10590.6666
$ws.Range("M132").Value = -8060.6666
$ws.Range("H135").Value = 74704.664
$ws.Range("J135").Value = 74704.664
$ws.Range("L135").Value = 74704.664
$ws.Range("N135").Value = -84844.664
$ws.Range("H141").Value = 92749.5
$ws.Range("J141").Value = 92749.5
$ws.Range("L141").Value = 92749.5
$ws.Range("N141").Value = -103109.5
